# Updates cryptos list figures (price + 1h volume change) to match the
# latest GitHub Actions scrape. Numeric-looking Price values need a
# leading quote so Excel stores them as text (matching the original
# inlineStr cells, e.g. '555.55' not the number 555.55) instead of
# silently parsing them into floats and dropping formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.814.67'
$ws.Range("E2").Value = '  +4.59%  '

# Row 3
$ws.Range("D3").Value = '3.337.67'
$ws.Range("E3").Value = '  +4.49%  '

# Row 5
$ws.Range("D5").Value = '''555.55'
$ws.Range("E5").Value = '  +3.37%  '

# Row 6
$ws.Range("D6").Value = '''151.96'
$ws.Range("E6").Value = '  +4.87%  '

# Row 7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("D8").Value = '''0.529'
$ws.Range("E8").Value = '  +2.06%  '

# Row 9
$ws.Range("D9").Value = '''7.49'
$ws.Range("E9").Value = '  +2.30%  '

# Row 10
$ws.Range("E10").Value = '  +4.29%  '

# Row 11
$ws.Range("D11").Value = '''0.437'
$ws.Range("E11").Value = '  +1.88%  '

# Row 12
$ws.Range("D12").Value = '3.913.04'
$ws.Range("E12").Value = '  +4.49%  '

# Row 13
$ws.Range("E13").Value = '  -0.16%  '

# Row 14
$ws.Range("E14").Value = '  +3.75%  '

# Row 15
$ws.Range("E15").Value = '  +2.52%  '

# Row 16
$ws.Range("D16").Value = '62.838.65'
$ws.Range("E16").Value = '  +4.61%  '

# Row 17
$ws.Range("D17").Value = '3.315.24'
$ws.Range("E17").Value = '  +2.73%  '

# Row 18
$ws.Range("D18").Value = '''6.49'
$ws.Range("E18").Value = '  +4.72%  '

# Row 19
$ws.Range("D19").Value = '''13.75'
$ws.Range("E19").Value = '  +4.76%  '

# Row 20
$ws.Range("E20").Value = '  +1.50%  '

# Row 21
$ws.Range("D21").Value = '''388.15'
$ws.Range("E21").Value = '  +1.76%  '

# Row 22
$ws.Range("E22").Value = '  +0.28%  '

# Row 23
$ws.Range("E23").Value = '  +1.69%  '

# Row 24
$ws.Range("D24").Value = '''70.64'
$ws.Range("E24").Value = '  +0.88%  '

# Row 25
$ws.Range("E25").Value = '  +4.86%  '

# Row 26
$ws.Range("E26").Value = '  -0.15%  '

# Row 27
$ws.Range("E27").Value = '  +6.60%  '

# Row 28
$ws.Range("E28").Value = '  +0.03%  '

# Row 29
$ws.Range("D29").Value = '''6.46'
$ws.Range("E29").Value = '  +4.54%  '

# Row 30
$ws.Range("E30").Value = '  +3.57%  '

# Row 31
$ws.Range("D31").Value = '''22.98'
$ws.Range("E31").Value = '  +2.43%  '

# Row 32
$ws.Range("E32").Value = '  +2.94%  '

# Row 33
$ws.Range("E33").Value = '  +6.03%  '

# Row 34
$ws.Range("E34").Value = '  +2.52%  '

# Row 35
$ws.Range("D35").Value = '''160.60'
$ws.Range("E35").Value = '  +2.47%  '

# Row 36
$ws.Range("E36").Value = '  +9.77%  '

# Row 37
$ws.Range("D37").Value = '''1.88'
$ws.Range("E37").Value = '  +11.50%  '

# Row 38
$ws.Range("D38").Value = '''27.18'
$ws.Range("E38").Value = '  +5.93%  '

# Row 39
$ws.Range("D39").Value = '2.844.45'
$ws.Range("E39").Value = '  +2.80%  '

# Row 40
$ws.Range("E40").Value = '  +3.35%  '

# Row 41
$ws.Range("D41").Value = '''0.0312'
$ws.Range("E41").Value = '  +8.56%  '

# Row 42
$ws.Range("E42").Value = '  +1.04%  '

# Row 43
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '''0.749'
$ws.Range("E43").Value = '  +2.71%  '

# Row 44
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '''40.67'
$ws.Range("E44").Value = '  +2.29%  '

# Row 45
$ws.Range("E45").Value = '  +2.41%  '

# Row 46
$ws.Range("D46").Value = '3.382.21'
$ws.Range("E46").Value = '  +4.51%  '

# Row 47
$ws.Range("D47").Value = '''21.93'
$ws.Range("E47").Value = '  +6.46%  '

# Row 48
$ws.Range("E48").Value = '  +3.38%  '

# Row 49
$ws.Range("E49").Value = '  +1.34%  '

# Row 50
$ws.Range("D50").Value = '''0.802'
$ws.Range("E50").Value = '  +1.00%  '

# Row 51
$ws.Range("D51").Value = '''282.08'
$ws.Range("E51").Value = '  +7.46%  '
